# ---------------------------------------------------------------------------
# Applies the two edits captured in the target diff:
#   1) Slide 10 - the "직사각형 21" rectangle that holds the hyperlinked video
#      URL: move/resize it and swap the displayed URL for a shortened
#      youtu.be link (keeping the original hyperlink target untouched), and
#      swap the run's lang/altLang.
#   2) Slide 8  - the "그림 1" picture: crop a sliver off the top/bottom of
#      the source image and reflow its position/size to match.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Slide 10: hyperlinked URL rectangle --------------------------------
$s10  = $p.Slides.Item(10)
$rect = $s10.Shapes.Item(4)

# Reposition / resize (EMU 4294989,3244334 / 3403817x369332).
$rect.Left   = 338.1881102362205
$rect.Top    = 255.45937347874016
$rect.Width  = 268.01708661417325
$rect.Height = 29.081259842519685

# Swap the displayed text for the shortened URL (hyperlink target/rId2 is
# left exactly as-is, only the visible run text changes).
$tr = $rect.TextFrame.TextRange
$tr.Text = "https://youtu.be/aKe38vRvUQE"

# Swap the run language tag (ko-KR -> en-US) to match the now-Latin text.
$tr.LanguageID = "en-US"

# --- 2) Slide 8: crop + reposition the UML picture --------------------------
$s8  = $p.Slides.Item(8)
$pic = $s8.Shapes.Item(1)

# Crop a thin sliver off the top/bottom of the source image
# (t="546" b="856" permille of the native image height).
$pic.PictureFormat.CropTop    = 2.61261
$pic.PictureFormat.CropBottom = 4.09596

# Reflow the frame to the cropped size/position (EMU 148997,1118139 /
# 11958283x4776123).
$pic.Left   = 11.732047644094488
$pic.Top    = 88.04244094488189
$pic.Width  = 941.5970866141732
$pic.Height = 376.0726776653543
